$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (reuses existing Company/Product/Price/Date strings)
$ws.Range("I1").Value = "Company"
$ws.Range("J1").Value = "Product"
$ws.Range("K1").Value = "Price"
$ws.Range("L1").Value = "Date"

# Row 2 - PrimeABGB / AMD Ryzen 5 2600
$ws.Range("I2").Value = "PrimeABGB"
$ws.Range("J2").Value = "AMD RYZEN 5 2600 6-CORE 3.4 GHZ (3.9 GHZ MAX BOOST) SOCKET AM4 PROCESSOR"
$ws.Range("K2").Value = 17344.0
$ws.Range("L2").Value = "29/09/2018"

# Row 3 - PrimeABGB / ASUS ROG STRIX motherboard
$ws.Range("I3").Value = "PrimeABGB"
$ws.Range("J3").Value = "ASUS ROG STRIX B350-F GAMING AM4 AMD MOTHERBOARD"
$ws.Range("K3").Value = 11275.0
$ws.Range("L3").Value = "29/09/2018"

# Row 4 - PrimeABGB / G.Skill TridentZ RAM
$ws.Range("I4").Value = "PrimeABGB"
$ws.Range("J4").Value = "G.SKILL TRIDENTZ RGB 16GB (2 X 8GB) DDR4 DESKTOP RAM F4-3200C14D-16GTZR"
$ws.Range("K4").Value = 22524.0
$ws.Range("L4").Value = "29/09/2018"

# Row 5 - PrimeABGB / Fractal Design Meshify case
$ws.Range("I5").Value = "PrimeABGB"
$ws.Range("J5").Value = "FRACTAL DESIGN MESHIFY C BLACK CASE FD-CA-MESH-C-BKO-TGL"
$ws.Range("K5").Value = 9189.0
$ws.Range("L5").Value = "29/09/2018"
